$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.831.43"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "1.986.94"
$ws.Range("E3").Value = "  -2.97%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.05"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.93"
$ws.Range("E7").Value = "  +10.43%  "
$ws.Range("E9").Value = "  -6.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.364"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0741"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  -2.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.939"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.71"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "2.276.35"
$ws.Range("E15").Value = "  -3.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.35"
$ws.Range("E16").Value = "  -1.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.90"
$ws.Range("E17").Value = "  +10.76%  "
$ws.Range("D18").Value = "1.988.64"
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("D19").Value = "35.723.94"
$ws.Range("E19").Value = "  -2.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.55"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "0.0₃0848"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.21"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.57"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("E25").Value = "  +14.68%  "
$ws.Range("E26").Value = "  -3.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.09"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.36"
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.02"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0988"
$ws.Range("E33").Value = "  +12.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0601"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.43"
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("E36").Value = "  +10.72%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  -2.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.67"
$ws.Range("E39").Value = "  +12.79%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0943"
$ws.Range("E42").Value = "  +5.02%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0214"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.51"
$ws.Range("E45").Value = "  +4.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.76"
$ws.Range("E46").Value = "  +5.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "93.18"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").Value = "1.359.18"
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.68"
$ws.Range("E51").Value = "  +3.55%  "
